$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data with 2020-04-01 figures (columns D and F for rows 24-30)
$ws.Range("D24").Value = 13155
$ws.Range("F24").Value = 727

$ws.Range("D25").Value = 9387
$ws.Range("F25").Value = 923

$ws.Range("D26").Value = 4032
$ws.Range("F26").Value = 509

$ws.Range("D27").Value = 239
$ws.Range("F27").Value = 59

$ws.Range("D28").Value = 2352
$ws.Range("F28").Value = 563

$ws.Range("D29").Value = 5102
$ws.Range("F29").Value = 1049

$ws.Range("D30").Value = 931
$ws.Range("F30").Value = 156

# Update the active cell selection to match the author's cursor position
$ws.Range("F31").Select()
